$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    # Force the written value to remain plain text (matches the source
    # inlineStr cells) instead of Excel auto-coercing numeric-looking
    # strings (e.g. "511.98") into real numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Rows where both Price (D) and Volume(1h) (E) changed
Set-TextValue $ws.Range("D2") "68.491.91"
Set-TextValue $ws.Range("E2") "  -0.79%  "
Set-TextValue $ws.Range("D3") "3.835.61"
Set-TextValue $ws.Range("E3") "  -2.38%  "
Set-TextValue $ws.Range("D5") "511.98"
Set-TextValue $ws.Range("E5") "  +5.21%  "
Set-TextValue $ws.Range("D6") "138.54"
Set-TextValue $ws.Range("E6") "  -6.09%  "
Set-TextValue $ws.Range("D7") "0.601"
Set-TextValue $ws.Range("E7") "  -3.55%  "
Set-TextValue $ws.Range("D9") "0.701"
Set-TextValue $ws.Range("E9") "  -5.12%  "
Set-TextValue $ws.Range("D12") "40.98"
Set-TextValue $ws.Range("E12") "  -4.65%  "
Set-TextValue $ws.Range("D13") "10.18"
Set-TextValue $ws.Range("E13") "  -3.06%  "
Set-TextValue $ws.Range("D14") "4.443.47"
Set-TextValue $ws.Range("E14") "  -2.57%  "
Set-TextValue $ws.Range("D15") "21.60"
Set-TextValue $ws.Range("E15") "  +8.14%  "
Set-TextValue $ws.Range("D16") "3.847.94"
Set-TextValue $ws.Range("E16") "  -2.40%  "
Set-TextValue $ws.Range("D17") "14.10"
Set-TextValue $ws.Range("E17") "  -0.83%  "
Set-TextValue $ws.Range("D20") "68.406.94"
Set-TextValue $ws.Range("E20") "  -0.97%  "
Set-TextValue $ws.Range("D21") "413.25"
Set-TextValue $ws.Range("E21") "  -5.81%  "
Set-TextValue $ws.Range("D22") "3.38"
Set-TextValue $ws.Range("E22") "  -2.94%  "
Set-TextValue $ws.Range("D23") "13.83"
Set-TextValue $ws.Range("E23") "  -5.90%  "
Set-TextValue $ws.Range("D25") "3.87"
Set-TextValue $ws.Range("E25") "  +4.31%  "
Set-TextValue $ws.Range("D26") "11.27"
Set-TextValue $ws.Range("E26") "  -6.52%  "
Set-TextValue $ws.Range("D27") "10.32"
Set-TextValue $ws.Range("E27") "  -7.33%  "
Set-TextValue $ws.Range("D28") "35.08"
Set-TextValue $ws.Range("E28") "  -5.68%  "
Set-TextValue $ws.Range("D29") "672.85"
Set-TextValue $ws.Range("E29") "  -5.63%  "
Set-TextValue $ws.Range("D30") "13.01"
Set-TextValue $ws.Range("E30") "  -2.88%  "
Set-TextValue $ws.Range("D32") "2.75"
Set-TextValue $ws.Range("E32") "  -4.90%  "
Set-TextValue $ws.Range("D33") "65.14"
Set-TextValue $ws.Range("E33") "  +6.97%  "
Set-TextValue $ws.Range("D34") "6.13"
Set-TextValue $ws.Range("E34") "  +1.58%  "
Set-TextValue $ws.Range("D35") "0.432"
Set-TextValue $ws.Range("E35") "  -7.63%  "
Set-TextValue $ws.Range("D36") "39.17"
Set-TextValue $ws.Range("E36") "  -4.43%  "
Set-TextValue $ws.Range("D39") "0.999"
Set-TextValue $ws.Range("E39") "  +0.11%  "
Set-TextValue $ws.Range("D43") "3.10"
Set-TextValue $ws.Range("E43") "  +2.94%  "
Set-TextValue $ws.Range("D44") "2.73"
Set-TextValue $ws.Range("E44") "  -8.16%  "
Set-TextValue $ws.Range("D45") "3.34"
Set-TextValue $ws.Range("E45") "  -2.41%  "
Set-TextValue $ws.Range("D46") "0.137"
Set-TextValue $ws.Range("E46") "  -3.98%  "
Set-TextValue $ws.Range("D47") "2.90"
Set-TextValue $ws.Range("E47") "  -2.04%  "
Set-TextValue $ws.Range("D50") "3.23"
Set-TextValue $ws.Range("E50") "  -4.22%  "

# Rows where only Volume(1h) (E) changed
Set-TextValue $ws.Range("E10") "  -7.39%  "
Set-TextValue $ws.Range("E11") "  -9.42%  "
Set-TextValue $ws.Range("E18") "  -2.15%  "
Set-TextValue $ws.Range("E19") "  +4.01%  "
Set-TextValue $ws.Range("E24") "  -4.37%  "
Set-TextValue $ws.Range("E31") "  -6.66%  "
Set-TextValue $ws.Range("E42") "  -4.24%  "

# Rows whose entire content (Coin, Link, Price, Volume) changed
# (coins were reordered/replaced between snapshots)
Set-TextValue $ws.Range("B37") "PEPE"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D37") "0.0₃0815"
Set-TextValue $ws.Range("E37") "  -11.46%  "
Set-TextValue $ws.Range("B38") "Kaspa"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D38") "0.148"
Set-TextValue $ws.Range("E38") "  -0.84%  "
Set-TextValue $ws.Range("B40") "ThetaToken"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D40") "3.34"
Set-TextValue $ws.Range("E40") "  +8.55%  "
Set-TextValue $ws.Range("B41") "FirstDigitalUSD"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D41") "1.00"
Set-TextValue $ws.Range("E41") "  -0.08%  "
Set-TextValue $ws.Range("B48") "Maker"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D48") "2.704.18"
Set-TextValue $ws.Range("E48") "  +11.35%  "
Set-TextValue $ws.Range("B49") "Monero"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D49") "143.48"
Set-TextValue $ws.Range("E49") "  +0.64%  "
Set-TextValue $ws.Range("B51") "BabyDogeCoin"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D51") "0.0₆0332"
Set-TextValue $ws.Range("E51") "  -11.09%  "
